# Applies the text edits made to slide "1067" (the single slide in this
# deck) as captured by the commit's canonical-XML diff:
#
#   1. In the pink-circle "Item #" label (TextBox 15), the placeholder
#      "{ITEM# }" (stray space before the closing brace) is tightened to
#      "{ITEM#}".
#   2. In the details textbox (TextBox 18), the "Dims:" line's
#      "Item Width(Inch)" gains a space -> "Item Width (Inch)".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Edit 1: TextBox 15 ------------------------------------------------
$shLabel = $s.Shapes.Item("TextBox 15")
# This shape auto-fits its height to its text (<a:spAutoFit/>); remember
# the height so the box keeps its original size once the character is
# removed, matching the source document (only the text changed there).
$origHeight = $shLabel.Height
$trLabel = $shLabel.TextFrame.TextRange
$foundItem = $trLabel.Find("{ITEM# }", 0, $true, $false)
if ($foundItem -ne $null) {
    $foundItem.Text = "{ITEM#}"
}
# Restore the original height. A minuscule epsilon compensates for the
# single-precision float round trip through the Height property so the
# restored value maps back to the identical EMU size.
[void]($shLabel.Height = $origHeight + 0.00001)

# --- Edit 2: TextBox 18 -------------------------------------------------
$shDetails = $s.Shapes.Item("TextBox 18")
$trDetails = $shDetails.TextFrame.TextRange
$foundInch = $trDetails.Find("(Inch)", 0, $true, $false)
if ($foundInch -ne $null) {
    [void]$foundInch.InsertBefore(" ")
}
